# Generate Report for Handoff
#
# This localization-status report moved from "In Translation" to
# "Ready for handoff" and the handoff/generation timestamps were
# refreshed. Update every cell that shares the old text so the
# workbook's shared strings collapse back down to one entry per
# distinct value (matching how the original report was produced),
# and touch up the status-column widths that Excel auto-sized for the
# new (longer) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime ------------
# zh-cn handoff timestamp: 2016-09-05 09:20:24 -> 2016-09-05 09:20:58
$wsZhCn.Range("H2").Value = "2016-09-05 09:20:58"

# de-de handoff timestamp (and mirrored Overview column):
# 2016-09-05 09:20:30 -> 2016-09-05 09:21:10
$wsOverview.Range("G2").Value = "2016-09-05 09:21:10"
$wsDeDe.Range("H2").Value     = "2016-09-05 09:21:10"

# --- Column width refresh for the Status columns -------------------------
# The longer "Ready for handoff" text made Excel widen the Status columns.
# ColumnWidth is quantized to whole pixels by the engine, so pick the input
# that lands on the closest achievable width to the recorded value
# (17.2159881591797 characters).
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # Overview!E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # Overview!F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33   # zh-cn!C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33   # de-de!C (Status)
